# Realestate Update resale numbers 2025-01-09 22:28
# Appends a new data row (row 16) to the CityResaleNum sheet with the
# 2025-01-09 22:28:50 resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column (A): must stay literal text "2025-01-09", not become a date
# serial. Pre-format the cell as Text, assign the value, then strip the
# formatting back off so the stored cell carries no explicit style - exactly
# like the other rows above it.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2025-01-09"
$ws.Range("A16").ClearFormats()

$ws.Range("B16").Value = "22:28:50"
$ws.Range("C16").Value = "Thursday"

# --- Week column (D): must stay literal text "01", not become the number 1.
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "01"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = 126373
$ws.Range("F16").Value = 143668
$ws.Range("G16").Value = 169623
$ws.Range("H16").Value = 159581
$ws.Range("I16").Value = -1
$ws.Range("J16").Value = 142870
$ws.Range("K16").Value = -1
$ws.Range("L16").Value = -1
$ws.Range("M16").Value = 192767
$ws.Range("N16").Value = 115328
$ws.Range("O16").Value = 45823
$ws.Range("P16").Value = 28481
$ws.Range("Q16").Value = 65024
$ws.Range("R16").Value = -1
$ws.Range("S16").Value = 48569
$ws.Range("T16").Value = -1
